$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 95

# Date column: write as literal text "2024-01-26" (not an Excel date serial).
# Forcing NumberFormat to "@" before the assignment stops Excel's automatic
# date-recognition; ClearFormats() afterwards drops the temporary text format
# so the cell is left with the same (default/unstyled) look as the rest of
# the data rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-26"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "23:00:38"
$ws.Cells.Item($row, 3).Value = "Friday"

# Week column: keep the leading zero ("03") as text instead of the number 3.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "03"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 136825
$ws.Cells.Item($row, 6).Value = 141618
$ws.Cells.Item($row, 7).Value = 171520
$ws.Cells.Item($row, 8).Value = 149199
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 122642
$ws.Cells.Item($row, 11).Value = 223965
$ws.Cells.Item($row, 12).Value = 256786
$ws.Cells.Item($row, 13).Value = 185376
$ws.Cells.Item($row, 14).Value = 110044
$ws.Cells.Item($row, 15).Value = 41408
$ws.Cells.Item($row, 16).Value = 30832
$ws.Cells.Item($row, 17).Value = 73580
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42570
$ws.Cells.Item($row, 20).Value = -1
